$d = $word.ActiveDocument

# The change targets <w:docDefaults> in styles.xml, which is not reachable
# through the normal Style/Font/ParagraphFormat object model (those only
# ever write an explicit override onto the "Normal" style, never touch
# docDefaults). Go through WordOpenXML (the flat-OPC serialization of the
# whole package) and surgically replace the docDefaults run/paragraph
# property blocks with literal string replacement.

$xml = $d.WordOpenXML

$oldRPr = '<w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b w:val="0"/><w:i w:val="0"/><w:smallCaps w:val="0"/><w:strike w:val="0"/><w:color w:val="000000"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:u w:val="none"/><w:shd w:val="clear" w:fill="auto"/><w:vertAlign w:val="baseline"/><w:lang w:val="en"/></w:rPr>'
$newRPr = '<w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en"/></w:rPr>'

$oldPPr = '<w:pPr><w:keepNext w:val="0"/><w:keepLines w:val="0"/><w:widowControl/><w:pBdr><w:top w:val="nil" w:sz="0" w:space="0"/><w:left w:val="nil" w:sz="0" w:space="0"/><w:bottom w:val="nil" w:sz="0" w:space="0"/><w:right w:val="nil" w:sz="0" w:space="0"/><w:between w:val="nil" w:sz="0" w:space="0"/></w:pBdr><w:shd w:val="clear" w:fill="auto"/><w:spacing w:before="0" w:after="0" w:line="276" w:lineRule="auto"/><w:ind w:left="0" w:right="0" w:firstLine="0"/><w:contextualSpacing w:val="0"/><w:jc w:val="left"/></w:pPr>'
$newPPr = '<w:pPr><w:spacing w:line="276" w:lineRule="auto"/></w:pPr>'

if (-not $xml.Contains($oldRPr)) {
    $oldRPrAlt = '<w:rPr><w:rFonts w:ascii="Arial" w:cs="Arial" w:eastAsia="Arial" w:hAnsi="Arial"/><w:b w:val="0"/><w:i w:val="0"/><w:smallCaps w:val="0"/><w:strike w:val="0"/><w:color w:val="000000"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:u w:val="none"/><w:shd w:fill="auto" w:val="clear"/><w:vertAlign w:val="baseline"/><w:lang w:val="en"/></w:rPr>'
    $newRPrAlt = '<w:rPr><w:rFonts w:ascii="Arial" w:cs="Arial" w:eastAsia="Arial" w:hAnsi="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en"/></w:rPr>'
    $xml = $xml.Replace($oldRPrAlt, $newRPrAlt)
} else {
    $xml = $xml.Replace($oldRPr, $newRPr)
}

if (-not $xml.Contains($oldPPr)) {
    $oldPPrAlt = '<w:pPr><w:keepNext w:val="0"/><w:keepLines w:val="0"/><w:widowControl w:val="1"/><w:pBdr><w:top w:space="0" w:sz="0" w:val="nil"/><w:left w:space="0" w:sz="0" w:val="nil"/><w:bottom w:space="0" w:sz="0" w:val="nil"/><w:right w:space="0" w:sz="0" w:val="nil"/><w:between w:space="0" w:sz="0" w:val="nil"/></w:pBdr><w:shd w:fill="auto" w:val="clear"/><w:spacing w:after="0" w:before="0" w:line="276" w:lineRule="auto"/><w:ind w:left="0" w:right="0" w:firstLine="0"/><w:contextualSpacing w:val="0"/><w:jc w:val="left"/></w:pPr>'
    $newPPrAlt = '<w:pPr><w:spacing w:line="276" w:lineRule="auto"/></w:pPr>'
    $xml = $xml.Replace($oldPPrAlt, $newPPrAlt)
} else {
    $xml = $xml.Replace($oldPPr, $newPPr)
}

$d.WordOpenXML = $xml

Write-Output "docDefaults trimmed"
